$d = $word.ActiveDocument

# --- 1. Paragraph 1: append the time portion to the existing date text ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1)
$r1.InsertAfter(", _HH_:_MM_:_SS_")

# --- 2. Paragraph 2: _N_ -> _n_ (x3) and drop the extra _N_ before -_CT_ ---
$d.Content.Find.Execute("_N_, _partie", $false, $false, $false, $false, $false, $true, 1, $false, "_n_, _partie", 2) | Out-Null
$d.Content.Find.Execute("_groupe_ _N_ (_UE__N_-_CT__N_)", $false, $false, $false, $false, $false, $true, 1, $false, "_groupe_ _n_ (_UE_-_CT__n_)", 2) | Out-Null

# --- 3. Font change: Times -> Times New Roman everywhere ---
# Setting .Name on the whole document content updates ascii+hAnsi for every run
# without disturbing paragraph marks (collapsed ranges) along the way.
$d.Content.Font.Name = "Times New Roman"

# The complex-script font (w:cs) only updates for the first run of whatever
# range it is applied to, so walk paragraph by paragraph for that piece.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    $pr.Font.NameBi = "Times New Roman"
    $pr.Font.NameOther = "Times New Roman"
}

# --- 4. Remove the last three (now-blank) paragraphs ---
$count = $d.Paragraphs.Count
$firstToDelete = $d.Paragraphs.Item($count - 2)
$lastToDelete = $d.Paragraphs.Item($count)
$delRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$delRange.Delete()
